# Dodanie podziału treningu na części
# Adds a "Trening" column that splits the training session into parts
# ("Duża Gra" / "Mała Gra") and refreshes the GPS sample rows with the
# per-part data, including converting the Timestamp column to real
# Excel date/time serial values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Trening" header (F1), matching the style of the other headers ---
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "Trening"

# --- Replace/extend the data rows (2-13) ---
# Columns: Timestamp(serial), Seconds, Velocity, Acceleration_SMA, Velocity_Bin, Trening
$data = @(
  @(45684.59259826389, 600.4,   13,    2.681206703186037,  "10-15", "Duża Gra"),
  @(45684.59329155093, 660.3,   12.01, 2.070611562047684,  "10-15", "Duża Gra"),
  @(45684.59408090277, 728.5,   13.16, 2.018800173486983,  "10-15", "Duża Gra"),
  @(45684.59259479167, 600.1,   9.67,  2.415275829178948,  "5-10",  "Duża Gra"),
  @(45684.59319085648, 651.6,   9.77,  1.75988280773163,   "5-10",  "Duża Gra"),
  @(45684.5932880787,  660,     9.58,  1.877202578953334,  "5-10",  "Duża Gra"),
  @(45684.59870127315, 1127.7,  12.05, 3.402332067489623,  "10-15", "Mała Gra"),
  @(45684.60333090278, 1527.7,  13.67, 3.264711362974984,  "10-15", "Mała Gra"),
  @(45684.60469895833, 1645.9,  13.52, 3.512729729924884,  "10-15", "Mała Gra"),
  @(45684.59869895833, 1127.5,  9.300000000000001, 3.255094221660067, "5-10", "Mała Gra"),
  @(45684.60142118055, 1362.7,  8.69,  3.03275578362601,   "5-10",  "Mała Gra"),
  @(45684.6040994213,  1594.1,  9.710000000000001, 3.23534676006862,  "5-10", "Mała Gra")
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $r++
}

# --- Apply the date/time number format to the Timestamp column ---
# (applied once to A2 and then to the whole A2:A13 range so the final
# style table matches the expected two numFmt entries / single cellXf)
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
